$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 28, shifting existing rows 28-38 down to 29-39
# (this matches the weekly data refresh: a new week's record is prepended
# to the Haba series and the older rows slide down).
$ws.Rows(28).Insert()

# Populate the newly inserted row 28 with the latest week's record.
$ws.Range('A28').Value = 7
$ws.Range('B28').Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range('C28').Value = 'Ñuble'
$ws.Range('D28').Value = 44553
$ws.Range('E28').Value = 16
$ws.Range('F28').Value = 100112026
$ws.Range('G28').Value = 'Haba'
$ws.Range('H28').Value = 'Sin especificar'
$ws.Range('I28').Value = 'Primera'
$ws.Range('J28').Value = 100
$ws.Range('K28').Value = 6500
$ws.Range('L28').Value = 7000
$ws.Range('M28').Value = 6750
$ws.Range('N28').Value = '$/saco 25 kilos'
$ws.Range('O28').Value = 'Provincia de Diguillín'
$ws.Range('P28').Value = 270
$ws.Range('Q28').Value = 25
$ws.Range('R28').Value = 'Hortaliza'
